
# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml  -> "Office Theme" (Office colour palette) - used by the notes master
#   ppt/theme/theme2.xml  -> "Integral"     (green/teal palette)    - used by the slide master (the
#                                                                     one actually driving every
#                                                                     slide's look)
# The authored edit swaps the two themes' contents, so the slides end up rendered with the
# (former theme1) "Office" colour palette while the notes master ends up with the (former
# theme2) "Integral" palette.
#
# The live/visible theme (the one bound to the slide master and therefore to every slide) is
# reachable from the object model via Slide.ThemeColorScheme - each of its 12 slots maps 1:1 onto
# the <a:clrScheme> children in document order:
#   1 dk1   2 lt1   3 dk2   4 lt2   5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
# ThemeColor.RGB uses the classic OLE "0xBBGGRR" byte order, so a target "RRGGBB" hex value has to
# be byte-reversed before assigning it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette = the "Office Theme" colours (currently sitting unused in theme1.xml) that must
# become the live theme's colours.
$officeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x44546A,  # 3  dk2
    0xE7E6E6,  # 4  lt2
    0x5B9BD5,  # 5  accent1
    0xED7D31,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0xFFC000,  # 8  accent4
    0x4472C4,  # 9  accent5
    0x70AD47,  # 10 accent6
    0x0563C1,  # 11 hlink
    0x954F72   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $rrggbb = $officeColors[$i - 1]
    $r = ($rrggbb -shr 16) -band 0xFF
    $g = ($rrggbb -shr 8) -band 0xFF
    $b = $rrggbb -band 0xFF
    $bgr = ($b -shl 16) -bor ($g -shl 8) -bor $r
    $tcs.Colors($i).RGB = $bgr
}
